# Fix bug in sales subtraction on the "Scenario Summary" sheet.
# Row 13 formulas were subtracting the wrong reference rows
# (C27/C28 instead of C26/C27). Update them accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario Summary")

$ws.Range("C13").Formula = "='Budget Revenues'!C7 - C26"
$ws.Range("D13").Formula = "='Optimistic (Raw)'!C7 - C26"
$ws.Range("E13").Formula = "='Pessimistic (Raw)'!C7 - C26"
$ws.Range("F13").Formula = "='Budget Revenues'!E7 - C27"
$ws.Range("G13").Formula = "='Optimistic (Raw)'!D7 - C27"
$ws.Range("H13").Formula = "='Pessimistic (Raw)'!D7 - C27"

# Update the selected cell to reflect the author's final cursor position.
$ws.Range("M15").Select()

$wb.Application.CalculateFull()
